$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New backup motion sensor button row (row 7)
# Shared-string pool order matters: "Backup Motion Sensor Button" must be
# interned before "Backup Motion Sensor Button Gnd", so write E7 first.
$ws.Range("E7").Value = "Backup Motion Sensor Button"
$ws.Range("E7").Interior.Color = $ws.Range("E11").Interior.Color

$ws.Range("B7").Value = "Backup Motion Sensor Button Gnd"
$ws.Range("B7").Interior.Color = $ws.Range("B12").Interior.Color

$ws.Range("C7").Value = "g"
$ws.Range("C7").Interior.Color = $ws.Range("C12").Interior.Color

$ws.Range("F7").Value = 16
$ws.Range("F7").Interior.Color = $ws.Range("F11").Interior.Color

# Update selection to match target state
$ws.Range("B7").Select()
